{"js": "// Replace each three-digit-by-one-digit multiplication equation with its\n// new value. Every source string is unique in the document, so a\n// matchCase/matchWholeWord exact search-and-replace is safe and precise.\nconst replacements = [\n  [\"420\u00d77=2940\", \"333\u00d76=1998\"],\n  [\"965\u00d75=4825\", \"576\u00d75=2880\"],\n  [\"330\u00d74=1320\", \"568\u00d74=2272\"],\n  [\"459\u00d76=2754\", \"647\u00d78=5176\"],\n  [\"258\u00d73=774\", \"173\u00d73=519\"],\n  [\"335\u00d75=1675\", \"362\u00d76=2172\"],\n  [\"470\u00d77=3290\", \"537\u00d74=2148\"],\n  [\"462\u00d77=3234\", \"524\u00d75=2620\"],\n  [\"354\u00d79=3186\", \"704\u00d74=2816\"],\n  [\"336\u00d73=1008\", \"967\u00d73=2901\"],\n  [\"656\u00d79=5904\", \"300\u00d72=600\"],\n  [\"511\u00d74=2044\", \"302\u00d75=1510\"],\n  [\"495\u00d75=2475\", \"992\u00d74=3968\"],\n  [\"139\u00d72=278\", \"234\u00d72=468\"],\n  [\"411\u00d79=3699\", \"877\u00d79=7893\"],\n  [\"744\u00d73=2232\", \"960\u00d76=5760\"],\n  [\"888\u00d78=7104\", \"403\u00d77=2821\"],\n  [\"263\u00d78=2104\", \"252\u00d78=2016\"],\n  [\"447\u00d78=3576\", \"172\u00d79=1548\"],\n  [\"390\u00d77=2730\", \"703\u00d72=1406\"],\n  [\"177\u00d77=1239\", \"321\u00d79=2889\"],\n  [\"860\u00d74=3440\", \"156\u00d76=936\"],\n  [\"909\u00d77=6363\", \"967\u00d72=1934\"],\n  [\"656\u00d76=3936\", \"783\u00d75=3915\"],\n  [\"332\u00d79=2988\", \"104\u00d72=208\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication equation with its\n# new value. Every source string is unique in the document, so an exact\n# Find/Replace (match case, whole string) is safe and precise.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @{Old=\"420\u00d77=2940\"; New=\"333\u00d76=1998\"},\n  @{Old=\"965\u00d75=4825\"; New=\"576\u00d75=2880\"},\n  @{Old=\"330\u00d74=1320\"; New=\"568\u00d74=2272\"},\n  @{Old=\"459\u00d76=2754\"; New=\"647\u00d78=5176\"},\n  @{Old=\"258\u00d73=774\"; New=\"173\u00d73=519\"},\n  @{Old=\"335\u00d75=1675\"; New=\"362\u00d76=2172\"},\n  @{Old=\"470\u00d77=3290\"; New=\"537\u00d74=2148\"},\n  @{Old=\"462\u00d77=3234\"; New=\"524\u00d75=2620\"},\n  @{Old=\"354\u00d79=3186\"; New=\"704\u00d74=2816\"},\n  @{Old=\"336\u00d73=1008\"; New=\"967\u00d73=2901\"},\n  @{Old=\"656\u00d79=5904\"; New=\"300\u00d72=600\"},\n  @{Old=\"511\u00d74=2044\"; New=\"302\u00d75=1510\"},\n  @{Old=\"495\u00d75=2475\"; New=\"992\u00d74=3968\"},\n  @{Old=\"139\u00d72=278\"; New=\"234\u00d72=468\"},\n  @{Old=\"411\u00d79=3699\"; New=\"877\u00d79=7893\"},\n  @{Old=\"744\u00d73=2232\"; New=\"960\u00d76=5760\"},\n  @{Old=\"888\u00d78=7104\"; New=\"403\u00d77=2821\"},\n  @{Old=\"263\u00d78=2104\"; New=\"252\u00d78=2016\"},\n  @{Old=\"447\u00d78=3576\"; New=\"172\u00d79=1548\"},\n  @{Old=\"390\u00d77=2730\"; New=\"703\u00d72=1406\"},\n  @{Old=\"177\u00d77=1239\"; New=\"321\u00d79=2889\"},\n  @{Old=\"860\u00d74=3440\"; New=\"156\u00d76=936\"},\n  @{Old=\"909\u00d77=6363\"; New=\"967\u00d72=1934\"},\n  @{Old=\"656\u00d76=3936\"; New=\"783\u00d75=3915\"},\n  @{Old=\"332\u00d79=2988\"; New=\"104\u00d72=208\"}\n)\n\nforeach ($p in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Execute($p.Old, $false, $false, $false, $false, $false, $true, 1, $false, $p.New, 2)\n}\n"}
